$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "x" marker from D21 and D22 (they no longer count toward the "Buy" SUMIF)
$ws.Range("D21").Clear()
$ws.Range("D22").Clear()

# Update the selection to G16, as in the edited workbook
$ws.Range("G16").Select()
